$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.126.33"
$ws.Range("E2").Value = "  +1.72%  "

# Row 3
$ws.Range("D3").Value = "3.774.13"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = "621.31"
$cD.Style = "Normal"
$ws.Range("E5").Value = "  +3.97%  "

# Row 6
$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value = "164.68"
$cD.Style = "Normal"
$ws.Range("E6").Value = "  +1.10%  "

# Row 7
$ws.Range("D7").Value = "3.772.05"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value = "0.520"
$cD.Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "

# Row 10
$ws.Range("E10").Value = "  +2.44%  "

# Row 11
$cD = $ws.Range("D11")
$cD.NumberFormat = "@"
$cD.Value = "0.451"
$cD.Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "

# Row 12
$cD = $ws.Range("D12")
$cD.NumberFormat = "@"
$cD.Value = "6.63"
$cD.Style = "Normal"
$ws.Range("E12").Value = "  +1.48%  "

# Row 13
$cD = $ws.Range("D13")
$cD.NumberFormat = "@"
$cD.Value = "0.0000247"
$cD.Style = "Normal"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14
$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = "35.69"
$cD.Style = "Normal"
$ws.Range("E14").Value = "  +1.09%  "

# Row 15
$ws.Range("D15").Value = "4.402.81"
$ws.Range("E15").Value = "  -0.15%  "

# Row 16
$ws.Range("D16").Value = "3.853.44"
$ws.Range("E16").Value = "  +2.25%  "

# Row 17
$ws.Range("D17").Value = "69.063.52"
$ws.Range("E17").Value = "  +1.63%  "

# Row 18
$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value = "17.67"
$cD.Style = "Normal"
$ws.Range("E18").Value = "  -2.76%  "

# Row 19
$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = "7.09"
$cD.Style = "Normal"
$ws.Range("E19").Value = "  +1.40%  "

# Row 20
$ws.Range("E20").Value = "  -1.34%  "

# Row 21
$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = "468.36"
$cD.Style = "Normal"
$ws.Range("E21").Value = "  +1.87%  "

# Row 22
$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = "9.59"
$cD.Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$cD = $ws.Range("D23")
$cD.NumberFormat = "@"
$cD.Value = "0.701"
$cD.Style = "Normal"
$ws.Range("E23").Value = "  +0.90%  "

# Row 24
$ws.Range("E24").Value = "  +4.32%  "

# Row 25
$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value = "83.17"
$cD.Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$cD = $ws.Range("D26")
$cD.NumberFormat = "@"
$cD.Value = "11.99"
$cD.Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "

# Row 27
$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value = "2.15"
$cD.Style = "Normal"
$ws.Range("E27").Value = "  +3.69%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cD = $ws.Range("D28")
$cD.NumberFormat = "@"
$cD.Value = "10.01"
$cD.Style = "Normal"
$ws.Range("E28").Value = "  +1.08%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value = "1.00"
$cD.Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").Value = "3.917.98"
$ws.Range("E30").Value = "  -0.27%  "

# Row 31
$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value = "2.25"
$cD.Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "

# Row 32
$ws.Range("E32").Value = "  +4.00%  "

# Row 33
$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = "7.29"
$cD.Style = "Normal"
$ws.Range("E33").Value = "  +1.01%  "

# Row 34
$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = "28.87"
$cD.Style = "Normal"
$ws.Range("E34").Value = "  -0.48%  "

# Row 35
$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value = "1.00"
$cD.Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cD = $ws.Range("D36")
$cD.NumberFormat = "@"
$cD.Value = "8.98"
$cD.Style = "Normal"
$ws.Range("E36").Value = "  +0.74%  "

# Row 37
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.722.57"
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = "0.158"
$cD.Style = "Normal"
$ws.Range("E38").Value = "  +12.72%  "

# Row 39
$ws.Range("E39").Value = "  +3.24%  "

# Row 40
$cD = $ws.Range("D40")
$cD.NumberFormat = "@"
$cD.Value = "3.36"
$cD.Style = "Normal"
$ws.Range("E40").Value = "  +4.93%  "

# Row 41
$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value = "5.79"
$cD.Style = "Normal"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42
$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value = "0.966"
$cD.Style = "Normal"
$ws.Range("E42").Value = "  -1.43%  "

# Row 43
$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value = "0.999"
$cD.Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = "0.300"
$cD.Style = "Normal"
$ws.Range("E45").Value = "  +1.62%  "

# Row 46
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cD = $ws.Range("D46")
$cD.NumberFormat = "@"
$cD.Value = "154.24"
$cD.Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "

# Row 47
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value = "42.92"
$cD.Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value = "46.79"
$cD.Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

# Row 49
$ws.Range("E49").Value = "  +3.24%  "

# Row 50
$cD = $ws.Range("D50")
$cD.NumberFormat = "@"
$cD.Value = "8.40"
$cD.Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "

# Row 51
$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = "1.37"
$cD.Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
